$d = $word.ActiveDocument
$d.Content.Find.Execute(
    "prototypes that test basic functionality",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "prototypes and the software sample available on the Collab site to test basic functionality",
    2
)
